# Auto-generated edit script for cryptos.xlsx update
# Applies the Tue Jul 18 14:28:01 UTC 2023 GitHub Actions cryptos list refresh:
#  - New row inserted for LidoStakedEther at row 8, shifting rows 8-50 down by one
#  - Cronos (previously row 51) drops off the bottom of the 50-row list
#  - Price (D) and Volume(1h) (E) values refreshed for every data row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.952.47'
$ws.Range("E2").Value = '  -1.08%  '

# Row 3
$ws.Range("D3").Value = '1.899.46'
$ws.Range("E3").Value = '  -0.76%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9982'
$ws.Range("E4").Value = '  -0.24%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7591'
$ws.Range("E5").Value = '  +2.51%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '240.63'
$ws.Range("E6").Value = '  -1.29%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9993'
$ws.Range("E7").Value = '  -0.09%  '

# Row 8
$ws.Range("B8").Value = 'LidoStakedEther'
$ws.Range("C8").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D8").Value = '1.896.51'
$ws.Range("E8").Value = '  -0.22%  '

# Row 9
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3054'
$ws.Range("E9").Value = '  -3.02%  '

# Row 10
$ws.Range("B10").Value = 'Solana'
$ws.Range("C10").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '25.54'
$ws.Range("E10").Value = '  -5.97%  '

# Row 11
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06852'
$ws.Range("E11").Value = '  -2.29%  '

# Row 12
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07990'
$ws.Range("E12").Value = '  +0.25%  '

# Row 13
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.7522'
$ws.Range("E13").Value = '  -3.87%  '

# Row 14
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.896.16'
$ws.Range("E14").Value = '  -0.71%  '

# Row 15
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.229'
$ws.Range("E15").Value = '  -1.45%  '

# Row 16
$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.46'
$ws.Range("E16").Value = '  -0.75%  '

# Row 17
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '29.948.92'
$ws.Range("E17").Value = '  -0.92%  '

# Row 18
$ws.Range("B18").Value = 'Avalanche'
$ws.Range("C18").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.97'
$ws.Range("E18").Value = '  -2.81%  '

# Row 19
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.007'
$ws.Range("E19").Value = '  +2.32%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '243.06'
$ws.Range("E20").Value = '  -1.06%  '

# Row 21
$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.000007723'
$ws.Range("E21").Value = '  -1.74%  '

# Row 22
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9992'
$ws.Range("E22").Value = '  -0.09%  '

# Row 23
$ws.Range("B23").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D23").Value = '2.146.05'
$ws.Range("E23").Value = '  -1.11%  '

# Row 24
$ws.Range("B24").Value = 'BinanceUSD'
$ws.Range("C24").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9974'
$ws.Range("E24").Value = '  -0.32%  '

# Row 25
$ws.Range("B25").Value = 'Chainlink'
$ws.Range("C25").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.994'
$ws.Range("E25").Value = '  +4.74%  '

# Row 26
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.267'
$ws.Range("E26").Value = '  -2.80%  '

# Row 27
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '165.75'
$ws.Range("E27").Value = '  +0.27%  '

# Row 28
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.80'
$ws.Range("E28").Value = '  -1.28%  '

# Row 29
$ws.Range("B29").Value = 'Stellar'
$ws.Range("C29").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1298'
$ws.Range("E29").Value = '  +1.69%  '

# Row 30
$ws.Range("B30").Value = 'LidoDAOToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.043'
$ws.Range("E30").Value = '  -4.02%  '

# Row 31
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.344'
$ws.Range("E31").Value = '  -1.10%  '

# Row 32
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.518'
$ws.Range("E32").Value = '  -2.02%  '

# Row 33
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.298'
$ws.Range("E33").Value = '  -0.85%  '

# Row 34
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.038'
$ws.Range("E34").Value = '  -1.22%  '

# Row 35
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.05381'
$ws.Range("E35").Value = '  +2.94%  '

# Row 36
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.255'
$ws.Range("E36").Value = '  -4.16%  '

# Row 37
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.7308'
$ws.Range("E37").Value = '  -3.06%  '

# Row 38
$ws.Range("B38").Value = 'HuobiToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.728'
$ws.Range("E38").Value = '  -1.12%  '

# Row 39
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01930'
$ws.Range("E39").Value = '  -1.15%  '

# Row 40
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.772'
$ws.Range("E40").Value = '  -0.74%  '

# Row 41
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.195'
$ws.Range("E41").Value = '  -3.48%  '

# Row 42
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4425'
$ws.Range("E42").Value = '  -1.69%  '

# Row 43
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '72.50'
$ws.Range("E43").Value = '  -4.72%  '

# Row 44
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.921'
$ws.Range("E44").Value = '  -1.31%  '

# Row 45
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9998'
$ws.Range("E45").Value = '  +0.05%  '

# Row 46
$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.8302'
$ws.Range("E46").Value = '  -0.24%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.617'
$ws.Range("E47").Value = '  -2.23%  '

# Row 48
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '101.11'
$ws.Range("E48").Value = '  -0.12%  '

# Row 49
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.766'
$ws.Range("E49").Value = '  -1.42%  '

# Row 50
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.049.59'
$ws.Range("E50").Value = '  -0.64%  '

# Row 51
$ws.Range("B51").Value = 'Elrond'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '36.28'
$ws.Range("E51").Value = '  -2.82%  '
